# Generate Report for Handback
# The handback transform failed because the handback file name did not match
# the expected handoff-derived file name. Update the status for the affected
# file (row 3, "38ee7d8b-50c3-4760-9f34-45bca1c8ac54") on the Overview sheet
# and on each locale sheet, and record the error detail per locale.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update status text everywhere it is used (shared string), by setting the
# cell values on each sheet where the status for this file is shown.
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handback transform failed"
$overview.Range("C3").Value = "Handback transform failed"

# Record the error detail explaining the handback file name mismatch.
$zhcn.Range("K3").Value = "Handback file name: fjogegqg.nau is different with handoff file name: 38ee7d8b-50c3-4760-9f34-45bca1c8ac54.319edca6af5cb0526285b7ee263e5d27ad37f7f0.zh-cn."
$dede.Range("K3").Value = "Handback file name: fjogegqg.nau is different with handoff file name: 38ee7d8b-50c3-4760-9f34-45bca1c8ac54.319edca6af5cb0526285b7ee263e5d27ad37f7f0.de-de."
